# Plotting, presentation, 2017 Snow Crab Data
# Applies the edits captured in the target diff: a handful of new "S" column
# zero entries, a couple of number-format tweaks (date -> general) on the
# "T" column, an updated "date disentangled" (P9), a restated "days to
# disentanglement" formula (U9) with an explanatory comment (V9), and a
# corrected max drag value (Z9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Conf gear free" (S) column entries -------------------------------
# Plain value, no special number format (matches existing A/B column style).
$ws.Range("S3").Value = 0
$ws.Range("S11").Value = 0

# Cells that pick up the worksheet's date display (same "d-mmm-yy" format
# used throughout columns M:Q) even though the stored value is 0.
$ws.Range("S4").Value = 0
$ws.Range("S4").NumberFormat = $ws.Range("M4").NumberFormat

$ws.Range("S6").Value = 0
$ws.Range("S6").NumberFormat = $ws.Range("M6").NumberFormat

$ws.Range("S7").Value = 0
$ws.Range("S7").NumberFormat = $ws.Range("M7").NumberFormat

$ws.Range("S9").Value = 0
$ws.Range("S9").NumberFormat = $ws.Range("M9").NumberFormat

$ws.Range("S10").Value = 0
$ws.Range("S10").NumberFormat = $ws.Range("N10").NumberFormat

# S8 already existed as a blank, styled cell -- just give it its value.
$ws.Range("S8").Value = 0

# --- Day-count ("T") column: switch from date display back to General -----
$ws.Range("T8").NumberFormat = "General"
$ws.Range("T9").NumberFormat = "General"
$ws.Range("T10").NumberFormat = "General"

# --- Row 9: updated disentanglement date & recomputed day counts ----------
# Date disentangled moves later (Amy's corrected report).
$ws.Range("P9").Value = 38245

# "Days to disentanglement" is now a live formula instead of a hard-coded
# number, and a comment cell explains the discrepancy with Amy's report.
$ws.Range("U9").Formula = "=P9-M9"
$ws.Range("V9").Value = "Amy's report says 808, but 433 days since LSGF"

# Corrected max drag value for this whale.
$ws.Range("Z9").Value = 24

# --- Selection / view bookkeeping ------------------------------------------
$ws.Range("G21:G23").Select()
